$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")

# ALC row 21
$ws_ALC.Range("H21").Value = 23836.166
$ws_ALC.Range("I21").Value = 18672.334
$ws_ALC.Range("K21").Value = 18672.334
$ws_ALC.Range("M21").Value = -18204.334

# ALC row 23
$ws_ALC.Range("H23").Value = 23836.166
$ws_ALC.Range("I23").Value = 18672.334
$ws_ALC.Range("K23").Value = 18672.334
$ws_ALC.Range("M23").Value = -18438.334

# ALC row 28
$ws_ALC.Range("H28").Value = 427894.94
$ws_ALC.Range("I28").Value = 653810.2
$ws_ALC.Range("J28").Value = 1166.1111
$ws_ALC.Range("K28").Value = 653810.2
$ws_ALC.Range("L28").Value = 1166.1111
$ws_ALC.Range("M28").Value = -653325.2
$ws_ALC.Range("N28").Value = -2136.1111

# ALC row 64
$ws_ALC.Range("H64").Value = 480812.47
$ws_ALC.Range("I64").Value = 1114531.4
$ws_ALC.Range("J64").Value = 5523.3335
$ws_ALC.Range("K64").Value = 1114531.4
$ws_ALC.Range("L64").Value = 5523.3335
$ws_ALC.Range("M64").Value = -1114283.4
$ws_ALC.Range("N64").Value = -6019.3335

# ALC row 67
$ws_ALC.Range("H67").Value = 480812.47
$ws_ALC.Range("I67").Value = 1114531.4
$ws_ALC.Range("J67").Value = 5523.3335
$ws_ALC.Range("K67").Value = 1114531.4
$ws_ALC.Range("L67").Value = 5523.3335
$ws_ALC.Range("M67").Value = -1113673.4
$ws_ALC.Range("N67").Value = -7239.3335

# ALC row 112
$ws_ALC.Range("H112").Value = 9405473
$ws_ALC.Range("J112").Value = 10102119
$ws_ALC.Range("L112").Value = 30306357
$ws_ALC.Range("N112").Value = -30308573

# ALC row 113
$ws_ALC.Range("H113").Value = 116810.555
$ws_ALC.Range("I113").Value = 149613.58
$ws_ALC.Range("J113").Value = 2000
$ws_ALC.Range("K113").Value = 149613.58
$ws_ALC.Range("L113").Value = 2000
$ws_ALC.Range("M113").Value = -146359.58
$ws_ALC.Range("N113").Value = -8508

# ALC row 116
$ws_ALC.Range("H116").Value = 15376666
$ws_ALC.Range("I116").Value = 15376666
$ws_ALC.Range("K116").Value = 15376666
$ws_ALC.Range("M116").Value = -15373224

# ALC row 129
$ws_ALC.Range("H129").Value = 1028.6451
$ws_ALC.Range("J129").Value = 1402
$ws_ALC.Range("L129").Value = 4206
$ws_ALC.Range("N129").Value = -14206

# ALC row 132
$ws_ALC.Range("H132").Value = 270815.6
$ws_ALC.Range("I132").Value = 406078.8
$ws_ALC.Range("K132").Value = 1218236.4
$ws_ALC.Range("M132").Value = -1215706.4

# ALC row 133
$ws_ALC.Range("H133").Value = 12830
$ws_ALC.Range("J133").Value = 12830
$ws_ALC.Range("L133").Value = 12830
$ws_ALC.Range("N133").Value = -22950

# ALC row 137
$ws_ALC.Range("H137").Value = 45456268
$ws_ALC.Range("I137").Value = 62501068
$ws_ALC.Range("J137").Value = 3467.1667
$ws_ALC.Range("K137").Value = 187503204
$ws_ALC.Range("L137").Value = 10401.5001
$ws_ALC.Range("M137").Value = -187500654
$ws_ALC.Range("N137").Value = -15501.5001

# ARM row 32
$ws_ARM.Range("H32").Value = 3926.8167
$ws_ARM.Range("I32").Value = 1909.8478
$ws_ARM.Range("K32").Value = 1909.8478
$ws_ARM.Range("M32").Value = -1622.8478

# ARM row 74
$ws_ARM.Range("H74").Value = 5911.852
$ws_ARM.Range("I74").Value = 1489.7894
$ws_ARM.Range("J74").Value = 16414.25
$ws_ARM.Range("K74").Value = 1489.7894
$ws_ARM.Range("L74").Value = 16414.25
$ws_ARM.Range("M74").Value = -615.7893999999999
$ws_ARM.Range("N74").Value = -18162.25

# ARM row 77
$ws_ARM.Range("H77").Value = 5911.852
$ws_ARM.Range("I77").Value = 1489.7894
$ws_ARM.Range("J77").Value = 16414.25
$ws_ARM.Range("K77").Value = 7448.946999999999
$ws_ARM.Range("L77").Value = 82071.25
$ws_ARM.Range("M77").Value = -3080.946999999999
$ws_ARM.Range("N77").Value = -90807.25

# ARM row 133
$ws_ARM.Range("H133").Value = 33252.2
$ws_ARM.Range("J133").Value = 33252.2
$ws_ARM.Range("L133").Value = 33252.2
$ws_ARM.Range("N133").Value = -38312.2

# ARM row 139
$ws_ARM.Range("H139").Value = 43959.285
$ws_ARM.Range("J139").Value = 43959.285
$ws_ARM.Range("L139").Value = 43959.285
$ws_ARM.Range("N139").Value = -54239.285

# BSM row 7
$ws_BSM.Range("H7").Value = 5001800
$ws_BSM.Range("I7").Value = 2400
$ws_BSM.Range("J7").Value = 20000000
$ws_BSM.Range("K7").Value = 2400
$ws_BSM.Range("L7").Value = 20000000
$ws_BSM.Range("M7").Value = -2287
$ws_BSM.Range("N7").Value = -20000226

# BSM row 107
$ws_BSM.Range("H107").Value = 816
$ws_BSM.Range("I107").Value = 700
$ws_BSM.Range("J107").Value = 839.2
$ws_BSM.Range("K107").Value = 700
$ws_BSM.Range("L107").Value = 839.2
$ws_BSM.Range("M107").Value = 1220
$ws_BSM.Range("N107").Value = -4679.2

# CRP row 10
$ws_CRP.Range("H10").Value = 466
$ws_CRP.Range("I10").Value = 466
$ws_CRP.Range("J10").Value = 0
$ws_CRP.Range("K10").Value = 466
$ws_CRP.Range("L10").Value = 0
$ws_CRP.Range("M10").ClearContents()
$ws_CRP.Range("N10").Value = -327

# CRP row 13
$ws_CRP.Range("H13").Value = 87783.336
$ws_CRP.Range("I13").Value = 200
$ws_CRP.Range("J13").Value = 105300
$ws_CRP.Range("K13").Value = 200
$ws_CRP.Range("L13").Value = 105300
$ws_CRP.Range("M13").Value = -61
$ws_CRP.Range("N13").Value = -105578

# CRP row 14
$ws_CRP.Range("H14").Value = 24000
$ws_CRP.Range("I14").Value = 1000
$ws_CRP.Range("J14").Value = 70000
$ws_CRP.Range("K14").Value = 1000
$ws_CRP.Range("L14").Value = 70000
$ws_CRP.Range("M14").Value = -830
$ws_CRP.Range("N14").Value = -70340

# CRP row 94
$ws_CRP.Range("H94").Value = 1451.8889
$ws_CRP.Range("I94").Value = 834.8333
$ws_CRP.Range("J94").Value = 1760.4166
$ws_CRP.Range("K94").Value = 834.8333
$ws_CRP.Range("L94").Value = 1760.4166
$ws_CRP.Range("M94").Value = -383.8333
$ws_CRP.Range("N94").Value = -2662.4166

# CRP row 99
$ws_CRP.Range("H99").Value = 15626004
$ws_CRP.Range("I99").Value = 31251000
$ws_CRP.Range("J99").Value = 1007
$ws_CRP.Range("K99").Value = 31251000
$ws_CRP.Range("L99").Value = 1007
$ws_CRP.Range("M99").Value = -31249502
$ws_CRP.Range("N99").Value = -4003

# CRP row 105
$ws_CRP.Range("H105").Value = 669.25
$ws_CRP.Range("I105").Value = 625.6667
$ws_CRP.Range("J105").Value = 800
$ws_CRP.Range("K105").Value = 625.6667
$ws_CRP.Range("L105").Value = 800
$ws_CRP.Range("M105").Value = 1121.3333
$ws_CRP.Range("N105").Value = -4294

# CRP row 126
$ws_CRP.Range("H126").Value = 15626004
$ws_CRP.Range("I126").Value = 31251000
$ws_CRP.Range("J126").Value = 1007
$ws_CRP.Range("K126").Value = 93753000
$ws_CRP.Range("L126").Value = 3021
$ws_CRP.Range("M126").Value = -93750530
$ws_CRP.Range("N126").Value = -7961

# CRP row 132
$ws_CRP.Range("H132").Value = 2178
$ws_CRP.Range("I132").Value = 1524.2903
$ws_CRP.Range("J132").Value = 4020.2727
$ws_CRP.Range("K132").Value = 4572.8709
$ws_CRP.Range("L132").Value = 12060.8181
$ws_CRP.Range("M132").Value = -2042.8709
$ws_CRP.Range("N132").Value = -17120.8181

# CUL row 50
$ws_CUL.Range("H50").Value = 336
$ws_CUL.Range("J50").Value = 349.5
$ws_CUL.Range("L50").Value = 1048.5
$ws_CUL.Range("N50").Value = -2010.5

# CUL row 53
$ws_CUL.Range("H53").Value = 336
$ws_CUL.Range("J53").Value = 349.5
$ws_CUL.Range("L53").Value = 1048.5
$ws_CUL.Range("N53").Value = -2010.5

# CUL row 56
$ws_CUL.Range("H56").Value = 4385.5
$ws_CUL.Range("I56").Value = 4385.5
$ws_CUL.Range("K56").Value = 4385.5
$ws_CUL.Range("M56").Value = -3855.5

# CUL row 131
$ws_CUL.Range("H131").Value = 1592.5714
$ws_CUL.Range("J131").Value = 1667.7885
$ws_CUL.Range("L131").Value = 5003.3655
$ws_CUL.Range("N131").Value = -15083.3655

# CUL row 138
$ws_CUL.Range("H138").Value = 819.2308
$ws_CUL.Range("I138").Value = 819.2308
$ws_CUL.Range("K138").Value = 2457.6924
$ws_CUL.Range("M138").Value = 2682.3076

# CUL row 140
$ws_CUL.Range("H140").Value = 7908.4707
$ws_CUL.Range("I140").Value = 21138
$ws_CUL.Range("J140").Value = 2396.1667
$ws_CUL.Range("K140").Value = 63414
$ws_CUL.Range("L140").Value = 7188.500100000001
$ws_CUL.Range("M140").Value = -58234
$ws_CUL.Range("N140").Value = -17548.5001

# GSM row 46
$ws_GSM.Range("H46").Value = 15510.6
$ws_GSM.Range("I46").Value = 7526.5
$ws_GSM.Range("J46").Value = 20833.334
$ws_GSM.Range("K46").Value = 7526.5
$ws_GSM.Range("L46").Value = 20833.334
$ws_GSM.Range("M46").Value = -7370.5
$ws_GSM.Range("N46").Value = -21145.334

# GSM row 132
$ws_GSM.Range("H132").Value = 3264.5454
$ws_GSM.Range("I132").Value = 2951.6191
$ws_GSM.Range("J132").Value = 3812.1667
$ws_GSM.Range("K132").Value = 8854.8573
$ws_GSM.Range("L132").Value = 11436.5001
$ws_GSM.Range("M132").Value = -6324.8573
$ws_GSM.Range("N132").Value = -16496.5001

# GSM row 138
$ws_GSM.Range("H138").Value = 58360
$ws_GSM.Range("J138").Value = 58360
$ws_GSM.Range("L138").Value = 58360
$ws_GSM.Range("N138").Value = -68640
